$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 0.1423950979931163
$arr[0,1] = 0
$arr[0,2] = 0.234083573227565
$arr[0,3] = 0.1782150073766182
$arr[0,4] = 1.105796097394865
$arr[0,5] = 0.5444902194853185
$arr[0,6] = 0.6826083526132365
$arr[0,7] = 0
$arr[0,8] = 0.1824366652562546
$arr[0,9] = 0.4539634355244289
$arr[0,10] = 0
$arr[0,11] = 0.2204584592205094
$arr[0,12] = 1.544647681924555
$arr[0,13] = 2.426055038172038
$arr[1,0] = 0.132950717731859
$arr[1,1] = 0
$arr[1,2] = 0.2313334024654097
$arr[1,3] = 0.1775680223802958
$arr[1,4] = 1.106988621844266
$arr[1,5] = 0.5449838828416276
$arr[1,6] = 0.6862965369911507
$arr[1,7] = 0
$arr[1,8] = 0.1827850043821933
$arr[1,9] = 0.396367287497327
$arr[1,10] = 0
$arr[1,11] = 0.204966476955903
$arr[1,12] = 1.558012874060694
$arr[1,13] = 2.434573712127047
$arr[2,0] = 0.1272202556083926
$arr[2,1] = 0
$arr[2,2] = 0.2297371686545375
$arr[2,3] = 0.177249435626031
$arr[2,4] = 1.108287195819784
$arr[2,5] = 0.5456301449589063
$arr[2,6] = 0.6888390851203212
$arr[2,7] = 0
$arr[2,8] = 0.1830853044207075
$arr[2,9] = 0.3608888485047146
$arr[2,10] = 0
$arr[2,11] = 0.1955045848880985
$arr[2,12] = 1.566740715599817
$arr[2,13] = 2.441104757837266
$arr[3,0] = 0.124902437986222
$arr[3,1] = 0
$arr[3,2] = 0.2291100225081522
$arr[3,3] = 0.1771394341199972
$arr[3,4] = 1.108958934191655
$arr[3,5] = 0.54597979033759
$arr[3,6] = 0.6899451741312674
$arr[3,7] = 0
$arr[3,8] = 0.1832294315790755
$arr[3,9] = 0.3464033750092597
$arr[3,10] = 0
$arr[3,11] = 0.1916617003657777
$arr[3,12] = 1.570428648277268
$arr[3,13] = 2.444093401687923
$arr[4,0] = 0.1245186220495498
$arr[4,1] = 0
$arr[4,2] = 0.22900729756212
$arr[4,3] = 0.1771223671481685
$arr[4,4] = 1.109079090259421
$arr[4,5] = 0.546043060297265
$arr[4,6] = 0.6901330682751308
$arr[4,7] = 0
$arr[4,8] = 0.1832546782207096
$arr[4,9] = 0.3439964279286016
$arr[4,10] = 0
$arr[4,11] = 0.1910243804323386
$arr[4,12] = 1.571048958307493
$arr[4,13] = 2.444609429800266
$arr[5,0] = 0.1271889260425922
$arr[5,1] = 0
$arr[5,2] = 0.2297286161358301
$arr[5,3] = 0.177247871776828
$arr[5,4] = 1.108295677712668
$arr[5,5] = 0.5456345110206655
$arr[5,6] = 0.6888537187740127
$arr[5,7] = 0
$arr[5,8] = 0.183087160065611
$arr[5,9] = 0.3606936032779799
$arr[5,10] = 0
$arr[5,11] = 0.1954527057494104
$arr[5,12] = 1.566789920701094
$arr[5,13] = 2.441143738747414
$arr[6,0] = 0.1391245824801928
$arr[6,1] = 0
$arr[6,2] = 0.2331162002418097
$arr[6,3] = 0.1779756261065337
$arr[6,4] = 1.106089815508859
$arr[6,5] = 0.5445892123556746
$arr[6,6] = 0.6838223883448791
$arr[6,7] = 0
$arr[6,8] = 0.1825388492224107
$arr[6,9] = 0.4341285552933414
$arr[6,10] = 0
$arr[6,11] = 0.2151065567486583
$arr[6,12] = 1.549147820538138
$arr[6,13] = 2.428722438928645
$arr[7,0] = 0.1630663841022084
$arr[7,1] = 0
$arr[7,2] = 0.2404883108123954
$arr[7,3] = 0.1800252272070573
$arr[7,4] = 1.106252811542987
$arr[7,5] = 0.5452628758298204
$arr[7,6] = 0.6761584757688581
$arr[7,7] = 0
$arr[7,8] = 0.1821485042665856
$arr[7,9] = 0.5771916145416185
$arr[7,10] = 0
$arr[7,11] = 0.2540361333873165
$arr[7,12] = 1.518684685716437
$arr[7,13] = 2.414678710244971
$arr[8,0] = 0.1809757690043199
$arr[8,1] = 0
$arr[8,2] = 0.2463445446899613
$arr[8,3] = 0.1819086508437735
$arr[8,4] = 1.109103843363435
$arr[8,5] = 0.5474202467438261
$arr[8,6] = 0.671866523113124
$arr[8,7] = 0
$arr[8,8] = 0.1822784055005258
$arr[8,9] = 0.6816857370781406
$arr[8,10] = 0
$arr[8,11] = 0.2828635959376129
$arr[8,12] = 1.498815360300178
$arr[8,13] = 2.410645252272758
$arr[9,0] = 0.1891910804661592
$arr[9,1] = 0
$arr[9,2] = 0.2491033146293802
$arr[9,3] = 0.182847077470047
$arr[9,4] = 1.110992943521595
$arr[9,5] = 0.5487631900564054
$arr[9,6] = 0.670203878771531
$arr[9,7] = 0
$arr[9,8] = 0.182427817991929
$arr[9,9] = 0.7290815945719942
$arr[9,10] = 0
$arr[9,11] = 0.2960248626984523
$arr[9,12] = 1.490320150469579
$arr[9,13] = 2.410174415884484
$arr[10,0] = 0.1923116509357072
$arr[10,1] = 0
$arr[10,2] = 0.2501615190869586
$arr[10,3] = 0.1832141365241249
$arr[10,4] = 1.111793332582181
$arr[10,5] = 0.5493237380366196
$arr[10,6] = 0.66961587980839
$arr[10,7] = 0
$arr[10,8] = 0.1824973666323189
$arr[10,9] = 0.7470082941291878
$arr[10,10] = 0
$arr[10,11] = 0.3010152768013867
$arr[10,12] = 1.487181277993244
$arr[10,13] = 2.41019216848747
$arr[11,0] = 0.1916391560665858
$arr[11,1] = 0
$arr[11,2] = 0.2499330162788453
$arr[11,3] = 0.1831345642464655
$arr[11,4] = 1.111617174807307
$arr[11,5] = 0.5492007008598847
$arr[11,6] = 0.6697406663571854
$arr[11,7] = 0
$arr[11,8] = 0.1824818115678895
$arr[11,9] = 0.7431484136357085
$arr[11,10] = 0
$arr[11,11] = 0.2999402163392233
$arr[11,12] = 1.487853818804012
$arr[11,13] = 2.410179627369018
$arr[12,0] = 0.1894476202363364
$arr[12,1] = 0
$arr[12,2] = 0.2491901034260593
$arr[12,3] = 0.1828770414590366
$arr[12,4] = 1.111057088633743
$arr[12,5] = 0.5488082643471301
$arr[12,6] = 0.6701546702078076
$arr[12,7] = 0
$arr[12,8] = 0.1824332799812041
$arr[12,9] = 0.730556863704777
$arr[12,10] = 0
$arr[12,11] = 0.2964352979800609
$arr[12,12] = 1.490060349098606
$arr[12,13] = 2.410171947849818
$arr[13,0] = 0.1881064867755811
$arr[13,1] = 0
$arr[13,2] = 0.2487368051878178
$arr[13,3] = 0.1827208233478999
$arr[13,4] = 1.110725089152737
$arr[13,5] = 0.5485746585957685
$arr[13,6] = 0.6704136763026298
$arr[13,7] = 0
$arr[13,8] = 0.1824052414268138
$arr[13,9] = 0.7228413978553192
$arr[13,10] = 0
$arr[13,11] = 0.2942892754021003
$arr[13,12] = 1.49142207885663
$arr[13,13] = 2.410192772356822
$arr[14,0] = 0.1804402349644363
$arr[14,1] = 0
$arr[14,2] = 0.2461661499953465
$arr[14,3] = 0.1818489616954864
$arr[14,4] = 1.10899229134737
$arr[14,5] = 0.5473397587294357
$arr[14,6] = 0.6719810067228167
$arr[14,7] = 0
$arr[14,8] = 0.1822704569285847
$arr[14,9] = 0.6785854153371815
$arr[14,10] = 0
$arr[14,11] = 0.2820044062523621
$arr[14,12] = 1.49938147191137
$arr[14,13] = 2.410703461789922
$arr[15,0] = 0.1757545721745828
$arr[15,1] = 0
$arr[15,2] = 0.2446133265518569
$arr[15,3] = 0.1813349837915936
$arr[15,4] = 1.108080853103175
$arr[15,5] = 0.5466747987316438
$arr[15,6] = 0.6730166898561407
$arr[15,7] = 0
$arr[15,8] = 0.1822108909371849
$arr[15,9] = 0.6513994163834127
$arr[15,10] = 0
$arr[15,11] = 0.2744799891922156
$arr[15,12] = 1.50440344224176
$arr[15,13] = 2.411366041014873
$arr[16,0] = 0.1730659383383966
$arr[16,1] = 0
$arr[16,2] = 0.2437291072411654
$arr[16,3] = 0.1810470449857284
$arr[16,4] = 1.107612361430071
$arr[16,5] = 0.5463263607762627
$arr[16,6] = 0.6736396672020248
$arr[16,7] = 0
$arr[16,8] = 0.182185132592231
$arr[16,9] = 0.6357497264991991
$arr[16,10] = 0
$arr[16,11] = 0.2701566335054224
$arr[16,12] = 1.507343102335462
$arr[16,13] = 2.411875531282135
$arr[17,0] = 0.1721567245867845
$arr[17,1] = 0
$arr[17,2] = 0.2434312616293965
$arr[17,3] = 0.1809508755474525
$arr[17,4] = 1.1074633159282
$arr[17,5] = 0.5462142297952539
$arr[17,6] = 0.6738552838697132
$arr[17,7] = 0
$arr[17,8] = 0.1821778721226579
$arr[17,9] = 0.6304488069919501
$arr[17,10] = 0
$arr[17,11] = 0.2686935999677189
$arr[17,12] = 1.508347208863007
$arr[17,13] = 2.412070090258339
$arr[18,0] = 0.1762527043987774
$arr[18,1] = 0
$arr[18,2] = 0.2447777044552595
$arr[18,3] = 0.1813889022854838
$arr[18,4] = 1.108172108904206
$arr[18,5] = 0.5467420628276614
$arr[18,6] = 0.672903616714791
$arr[18,7] = 0
$arr[18,8] = 0.182216352001646
$arr[18,9] = 0.6542947680327416
$arr[18,10] = 0
$arr[18,11] = 0.275280513735737
$arr[18,12] = 1.503863551014803
$arr[18,13] = 2.41128222088625
$arr[19,0] = 0.1900910682321211
$arr[19,1] = 0
$arr[19,2] = 0.249407948901009
$arr[19,3] = 0.1829523650320546
$arr[19,4] = 1.111219292872178
$arr[19,5] = 0.5489221209630983
$arr[19,6] = 0.6700319385002018
$arr[19,7] = 0
$arr[19,8] = 0.1824471830473939
$arr[19,9] = 0.7342558880064303
$arr[19,10] = 0
$arr[19,11] = 0.297464602716623
$arr[19,12] = 1.489410119145973
$arr[19,13] = 2.410168883578905
$arr[20,0] = 0.1991911606250198
$arr[20,1] = 0
$arr[20,2] = 0.2525128264748417
$arr[20,3] = 0.1840423385522563
$arr[20,4] = 1.113706353366283
$arr[20,5] = 0.5506500431687584
$arr[20,6] = 0.668397627243948
$arr[20,7] = 0
$arr[20,8] = 0.1826736312431763
$arr[20,9] = 0.786391709196721
$arr[20,10] = 0
$arr[20,11] = 0.312001096377486
$arr[20,12] = 1.480419041160346
$arr[20,13] = 2.410583926937051
$arr[21,0] = 0.1943292204121576
$arr[21,1] = 0
$arr[21,2] = 0.2508485245306957
$arr[21,3] = 0.1834543769456154
$arr[21,4] = 1.112333660369003
$arr[21,5] = 0.5497000776651504
$arr[21,6] = 0.6692477221106259
$arr[21,7] = 0
$arr[21,8] = 0.1825458615275792
$arr[21,9] = 0.7585774903481877
$arr[21,10] = 0
$arr[21,11] = 0.3042393279002553
$arr[21,12] = 1.485176130978608
$arr[21,13] = 2.41025788727481
$arr[22,0] = 0.1760274825579558
$arr[22,1] = 0
$arr[22,2] = 0.2447033626601609
$arr[22,3] = 0.181364502201852
$arr[22,4] = 1.108130679249882
$arr[22,5] = 0.5467115472710447
$arr[22,6] = 0.6729546512623727
$arr[22,7] = 0
$arr[22,8] = 0.1822138566173024
$arr[22,9] = 0.6529858420691426
$arr[22,10] = 0
$arr[22,11] = 0.2749185886802152
$arr[22,12] = 1.504107472416969
$arr[22,13] = 2.411319715489554
$arr[23,0] = 0.156532873830173
$arr[23,1] = 0
$arr[23,2] = 0.2384163740229894
$arr[23,3] = 0.1794043110160679
$arr[23,4] = 1.10572870152113
$arr[23,5] = 0.5447888418653832
$arr[23,6] = 0.6779963861819311
$arr[23,7] = 0
$arr[23,8] = 0.1821808671797811
$arr[23,9] = 0.5385942842512463
$arr[23,10] = 0
$arr[23,11] = 0.2434641924751659
$arr[23,12] = 1.526484107527402
$arr[23,13] = 2.417374037348168

$ws.Range("B2:O25").Value = $arr
